# Updated legacy GSC export data:
# The first data row of the "Chart" sheet (2025-11-07, a now-obsolete/duplicate
# leading entry) is removed. Deleting the row shifts every subsequent row up
# by one and lets Excel recompute the shared-strings table / sheet dimension
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
